# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right before the existing "2021-Q3"
#    sheet (so tab order becomes: 总计, 2022-Q3, 2021-Q3, 2021-Q1, 2020-Q4).
# 2) Populate the new sheet with the fund-holding table for the new quarter.
# 3) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 right
#    after the header, shifting the existing quarters down, and append the
#    2020-Q4 row that falls off the bottom.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q3_21   = $wb.Worksheets.Item("2021-Q3")

# --- 1) Create the new sheet (positioned later, once it is populated) -----
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# --- 2) Fill in the "2022-Q3" sheet ----------------------------------------
# Header row (bold / centered / bordered, matching the other quarter sheets)
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Data row
$newSheet.Cells.Item(2,1).Value = 0

$textCells = $newSheet.Range("B2:G2")
$textCells.NumberFormat = "@"

$newSheet.Cells.Item(2,2).Value = "001978"
$newSheet.Cells.Item(2,3).Value = "泰信互联网+主题灵活配置混合"
$newSheet.Cells.Item(2,4).Value = "0.10"
$newSheet.Cells.Item(2,5).Value = "36.51"
$newSheet.Cells.Item(2,6).Value = "1.14"
$newSheet.Cells.Item(2,7).Value = "0.0011"
$newSheet.Cells.Item(2,8).Value = 6

# Copy the bold/centered/bordered style used by the header + index column of
# the sister sheets onto the matching cells of the new sheet.
$q3_21.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)

$q3_21.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Move the finished sheet into place, right before "2021-Q3", so the tab
# order becomes: 总计, 2022-Q3, 2021-Q3, 2021-Q1, 2020-Q4.
$newSheet.Move($q3_21)

# --- 3) Update the summary sheet -------------------------------------------
# Current rows 2-4 (2021-Q3 / 2021-Q1 / 2020-Q4) shift down by one row to
# rows 3-5, and a brand-new row 2 is inserted for 2022-Q3. The quarter labels
# in column B stay attached to the same row they were already on (Insert
# pushes them down for us); only the new row 2 needs its text filled in.
$summary.Rows("2:2").Insert(-4121)
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 1
$summary.Cells.Item(2,4).Value = 0

# Renumber the index column (A) for the rows that shifted down.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3

# Copy the bold/border style of the index column onto the new row's A cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
